$d = $word.ActiveDocument

# 1. "Expected Check-out Date" -> "Check-out Date" (table header cell)
$found = $d.Content.Find.Execute("Expected Check-out Date", $true, $false, $false, $false, $false, $true, 1, $false, "Check-out Date", 2)

# 2. Turn the trailing empty paragraph (after the table, before the sectPr)
#    into one that contains a "_GoBack" bookmark, matching Word's habit of
#    recording the last edit position.
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Start = $r.Start - 1
$bm = $d.Bookmarks.Add("_GoBack", $r)
